$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("管理员任务表")

# Insert a new row above row 7, shifting existing rows 7-13 down to 8-14
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new field documentation
$ws.Cells.Item(7, 2).Value = "name"
$ws.Cells.Item(7, 3).Value = "商户名"

# The row insert carries the old row 7's formatting onto the rows below it
# (now rows 8-11); clear that inherited formatting so it matches the plain,
# unformatted cells used for the rest of the "admin" task rows.
$ws.Range("D8:D11").ClearFormats()

# Update the selection to reflect where the user left off editing
$ws.Range("C7").Select()
